$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 21, shifting existing rows 21-23 down to 22-24.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly record.
$ws.Cells.Item(21, 1).Value = 8
$ws.Cells.Item(21, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(21, 3).Value = "Coquimbo"
$ws.Cells.Item(21, 4).Value = 44637
$ws.Cells.Item(21, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100101
$ws.Cells.Item(21, 8).Value = "Berries"
$ws.Cells.Item(21, 9).Value = 100101001
$ws.Cells.Item(21, 10).Value = "Arándano (blue)"
$ws.Cells.Item(21, 11).Value = "Sin especificar"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 160
$ws.Cells.Item(21, 14).Value = 6000
$ws.Cells.Item(21, 15).Value = 6500
$ws.Cells.Item(21, 16).Value = 6250
$ws.Cells.Item(21, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(21, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(21, 19).Value = 3125
$ws.Cells.Item(21, 20).Value = 2
